# The presentation previously used the "Integral" (Red Violet) design theme.
# This change swaps it for the default "Office Theme" colour scheme
# (the theme that used to live alongside it as the secondary/notes theme),
# i.e. applying a different Design to the deck.
#
# PowerPoint's Theme Colors are exposed on the (single) slide master's
# Theme object as a 12-slot ThemeColorScheme:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hlink, 12 folHlink
# Colours are RGB() encoded integers (R + G*256 + B*65536).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = RGB(0x00, 0x00, 0x00)   # dk1      -> 000000
$tcs.Item(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = RGB(0x44, 0x54, 0x6A)   # dk2      -> 44546A
$tcs.Item(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = RGB(0xED, 0x7D, 0x31)   # accent2  -> ED7D31
$tcs.Item(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # accent4  -> FFC000
$tcs.Item(9).RGB  = RGB(0x44, 0x72, 0xC4)   # accent5  -> 4472C4
$tcs.Item(10).RGB = RGB(0x70, 0xAD, 0x47)   # accent6  -> 70AD47
$tcs.Item(11).RGB = RGB(0x05, 0x63, 0xC1)   # hlink    -> 0563C1
$tcs.Item(12).RGB = RGB(0x95, 0x4F, 0x72)   # folHlink -> 954F72
